# C25 RRC pulse table scale
# Update the pulse-shaping table on Sheet1: bump the BPSK/IL2P payload value
# (D4) and the loop-gain Integral Limit (C11). Dependent formulas
# (E4, F4, F5, C12) recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = 75
$ws.Range("C11").Value = 30000

$excel.Calculate()
